$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.048.19"
$ws.Range("E2").Value = "  +2.88%  "
$ws.Range("D3").Value = "2.414.01"
$ws.Range("E3").Value = "  +3.77%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "559.01"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").Value = "138.22"
$ws.Range("E6").Value = "  +5.52%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "2.410.99"
$ws.Range("E9").Value = "  +3.78%  "
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("D11").Value = "5.74"
$ws.Range("E11").Value = "  +3.82%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "0.346"
$ws.Range("E13").Value = "  +3.36%  "
$ws.Range("D14").Value = "25.79"
$ws.Range("E14").Value = "  +9.29%  "
$ws.Range("D15").Value = "2.843.52"
$ws.Range("E15").Value = "  +3.73%  "
$ws.Range("D16").Value = "61.980.51"
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("E17").Value = "  +5.03%  "
$ws.Range("D18").Value = "2.419.24"
$ws.Range("E18").Value = "  +3.42%  "
$ws.Range("E19").Value = "  +4.66%  "
$ws.Range("D20").Value = "343.16"
$ws.Range("E20").Value = "  +9.34%  "
$ws.Range("D21").Value = "4.22"
$ws.Range("E21").Value = "  +2.08%  "
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "64.94"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "8.32"
$ws.Range("E27").Value = "  +5.94%  "
$ws.Range("D28").Value = "1.51"
$ws.Range("E28").Value = "  +12.04%  "
$ws.Range("E29").Value = "  +15.56%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0778"
$ws.Range("E30").Value = "  +6.40%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.79"
$ws.Range("E31").Value = "  +3.49%  "
$ws.Range("D32").Value = "6.34"
$ws.Range("E32").Value = "  +6.76%  "
$ws.Range("D33").Value = "171.35"
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "1.42"
$ws.Range("E34").Value = "  +3.18%  "
$ws.Range("B35").Value = "PolygonEcosystemToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D35").Value = "0.395"
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("D36").Value = "379.63"
$ws.Range("E36").Value = "  +17.52%  "
$ws.Range("D37").Value = "18.51"
$ws.Range("E37").Value = "  +3.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.50"
$ws.Range("E38").Value = "  +10.79%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").Value = "1.68"
$ws.Range("E41").Value = "  +9.42%  "
$ws.Range("D42").Value = "39.04"
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("D43").Value = "145.16"
$ws.Range("E43").Value = "  +5.28%  "
$ws.Range("E44").Value = "  +4.66%  "
$ws.Range("D45").Value = "20.62"
$ws.Range("E45").Value = "  +8.04%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.0956"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0518"
$ws.Range("E47").Value = "  +4.78%  "
$ws.Range("D48").Value = "0.586"
$ws.Range("E48").Value = "  +4.29%  "
$ws.Range("D49").Value = "17.89"
$ws.Range("E49").Value = "  +5.65%  "
$ws.Range("D50").Value = "0.0221"
$ws.Range("E50").Value = "  +3.46%  "
$ws.Range("D51").Value = "0.0₆0223"
$ws.Range("E51").Value = "  +3.37%  "
